$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule" (sheet1): new optimisation run (run 67) ----
$ws1 = $wb.Worksheets.Item("Schedule")

# Row 2
$ws1.Range("A2").Value2 = 46039.29166666666
$ws1.Range("B2").Value2 = 46039.95833333334
$ws1.Range("C2").Value2 = 16
$ws1.Range("D2").Value2 = 60.48
$ws1.Range("E2").Value2 = 785.4514005
$ws1.Range("F2").Value2 = 12.98696098710317

# Row 3
$ws1.Range("A3").Value2 = 46040.29166666666
$ws1.Range("B3").Value2 = 46040.79166666666
$ws1.Range("C3").Value2 = 12
$ws1.Range("D3").Value2 = 45.36
$ws1.Range("E3").Value2 = -103.04050575
$ws1.Range("F3").Value2 = -2.271616087962963

# Old rows 4 and 5 are no longer part of the result set -> remove them
$ws1.Rows("4:5").Delete()

# ---- Sheet "Detailed" (sheet2): updated price/status values ----
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("E3").Value2 = "OFF"
$ws2.Range("E4").Value2 = "OFF"
$ws2.Range("E5").Value2 = "OFF"
$ws2.Range("E6").Value2 = "OFF"
$ws2.Range("E7").Value2 = "OFF"
$ws2.Range("E8").Value2 = "OFF"
$ws2.Range("E9").Value2 = "OFF"
$ws2.Range("E10").Value2 = "OFF"
$ws2.Range("E26").Value2 = "ON"
$ws2.Range("E27").Value2 = "ON"
$ws2.Range("B42").Value2 = 9.637499999999999
$ws2.Range("E42").Value2 = "ON"
$ws2.Range("B43").Value2 = 16.85505
$ws2.Range("E43").Value2 = "ON"
$ws2.Range("B44").Value2 = 9.29172
$ws2.Range("E44").Value2 = "ON"
$ws2.Range("B45").Value2 = 36.2
$ws2.Range("C45").Value2 = "historical"
$ws2.Range("E45").Value2 = "ON"
$ws2.Range("C46").Value2 = "historical"
$ws2.Range("E46").Value2 = "ON"
$ws2.Range("B47").Value2 = 36.2
$ws2.Range("E47").Value2 = "ON"
$ws2.Range("B49").Value2 = 56.98
$ws2.Range("B51").Value2 = 36.2
$ws2.Range("B53").Value2 = 56.98
$ws2.Range("B55").Value2 = 55.10665
$ws2.Range("B64").Value2 = 8.13148
$ws2.Range("E64").Value2 = "ON"
$ws2.Range("B65").Value2 = 0.51
$ws2.Range("B66").Value2 = -5.51
$ws2.Range("B67").Value2 = 0.00962
$ws2.Range("B68").Value2 = 0
$ws2.Range("B69").Value2 = -0.32596
$ws2.Range("B70").Value2 = -5.50985
$ws2.Range("B71").Value2 = 0.6369899999999999
$ws2.Range("B72").Value2 = 0.62188
$ws2.Range("B73").Value2 = 0
$ws2.Range("B74").Value2 = 0
$ws2.Range("B75").Value2 = -5.51
$ws2.Range("B76").Value2 = -6.22899
$ws2.Range("B77").Value2 = -6.8
$ws2.Range("B79").Value2 = -23.5
$ws2.Range("B80").Value2 = -14.76423
$ws2.Range("B81").Value2 = -14.36017
$ws2.Range("B82").Value2 = -7.41863
$ws2.Range("B83").Value2 = -7.24787
$ws2.Range("B84").Value2 = -6.37751
$ws2.Range("B85").Value2 = -0.63398
$ws2.Range("B86").Value2 = -1.41571
$ws2.Range("B87").Value2 = 0.00036
$ws2.Range("B88").Value2 = 36.0601
$ws2.Range("E88").Value2 = "OFF"
$ws2.Range("B89").Value2 = 46.12851
$ws2.Range("B90").Value2 = 56.98
$ws2.Range("B91").Value2 = 55.8977
$ws2.Range("B92").Value2 = 44.47659
$ws2.Range("B94").Value2 = 54.94853
